$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values for columns P (14) and Q (15), matching existing header style
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update existing data values in columns B..H, J, K for rows 2-25
$ws.Range("B2").Value = 25.09703835460632
$ws.Range("C2").Value = 20.34457718097449
$ws.Range("D2").Value = 3.256916493061351
$ws.Range("E2").Value = 29.94318199716722
$ws.Range("F2").Value = 22.88231735487936
$ws.Range("G2").Value = 31.9857944657687
$ws.Range("H2").Value = 3.595090717041712
$ws.Range("J2").Value = 9.826643705189863
$ws.Range("K2").Value = 12.38257868810889
$ws.Range("B3").Value = 23.44037060755681
$ws.Range("C3").Value = 19.01244184219233
$ws.Range("D3").Value = 3.238608711798094
$ws.Range("E3").Value = 27.96966651360496
$ws.Range("F3").Value = 21.81932269912853
$ws.Range("G3").Value = 30.25526731914831
$ws.Range("H3").Value = 3.285551972856177
$ws.Range("J3").Value = 9.636311993851956
$ws.Range("K3").Value = 12.4078316464304
$ws.Range("B4").Value = 22.36475999134314
$ws.Range("C4").Value = 18.14868424621398
$ws.Range("D4").Value = 3.226650033764871
$ws.Range("E4").Value = 26.69359321050832
$ws.Range("F4").Value = 21.16699812314214
$ws.Range("G4").Value = 29.18497104870739
$ws.Range("H4").Value = 3.089952386807184
$ws.Range("J4").Value = 9.527406413215061
$ws.Range("K4").Value = 12.44280991433582
$ws.Range("B5").Value = 21.91167502323856
$ws.Range("C5").Value = 17.78512890841431
$ws.Range("D5").Value = 3.221601125125949
$ws.Range("E5").Value = 26.15731293373871
$ws.Range("F5").Value = 20.90156385787707
$ws.Range("G5").Value = 28.74733997123269
$ws.Range("H5").Value = 3.008804726189135
$ws.Range("J5").Value = 9.484982209435623
$ws.Range("K5").Value = 12.46161856841449
$ws.Range("B6").Value = 21.83554680979005
$ws.Range("C6").Value = 17.72406116072703
$ws.Range("D6").Value = 3.220752298085092
$ws.Range("E6").Value = 26.06727895313499
$ws.Range("F6").Value = 20.85752210515435
$ws.Range("G6").Value = 28.67459801585307
$ws.Range("H6").Value = 2.99524200038953
$ws.Range("J6").Value = 9.478054376997786
$ws.Range("K6").Value = 12.46500694907239
$ws.Range("B7").Value = 22.3587094058925
$ws.Range("C7").Value = 18.14382809604593
$ws.Range("D7").Value = 3.22658264694329
$ws.Range("E7").Value = 26.68642673509998
$ws.Range("F7").Value = 21.16341636066593
$ws.Range("G7").Value = 29.17907429254616
$ws.Range("H7").Value = 3.08886388729853
$ws.Range("J7").Value = 9.526826411419101
$ws.Range("K7").Value = 12.44304557986732
$ws.Range("B8").Value = 24.53789224471706
$ws.Range("C8").Value = 19.89473562398932
$ws.Range("D8").Value = 3.250755654125689
$ws.Range("E8").Value = 29.27600187156317
$ws.Range("F8").Value = 22.51589094607733
$ws.Range("G8").Value = 31.39097016296273
$ws.Range("H8").Value = 3.489483529931533
$ws.Range("J8").Value = 9.759339665456533
$ws.Range("K8").Value = 12.38703715021005
$ws.Range("B9").Value = 28.35186414356699
$ws.Range("C9").Value = 22.96743012778047
$ws.Range("D9").Value = 3.292277094527144
$ws.Range("E9").Value = 33.84901422894479
$ws.Range("F9").Value = 25.16000261213185
$ws.Range("G9").Value = 35.6512576612008
$ws.Range("H9").Value = 4.233339730280886
$ws.Range("J9").Value = 10.28025617375507
$ws.Range("K9").Value = 12.44714190048514
$ws.Range("B10").Value = 30.87152664672094
$ws.Range("C10").Value = 24.95319420236088
$ws.Range("D10").Value = 3.356618153387374
$ws.Range("E10").Value = 36.02948667031383
$ws.Range("F10").Value = 26.82482322773726
$ws.Range("G10").Value = 38.2560950477108
$ws.Range("H10").Value = 4.715839863095707
$ws.Range("J10").Value = 10.60390845687684
$ws.Range("K10").Value = 12.44557716714615
$ws.Range("B11").Value = 31.88981052527685
$ws.Range("C11").Value = 25.33437077529264
$ws.Range("D11").Value = 3.746802016271704
$ws.Range("E11").Value = 29.51772418503302
$ws.Range("F11").Value = 25.33217661157738
$ws.Range("G11").Value = 35.4651073883489
$ws.Range("H11").Value = 5.02979180152991
$ws.Range("J11").Value = 9.888245507499564
$ws.Range("K11").Value = 11.0078975655774
$ws.Range("B12").Value = 32.24315554862397
$ws.Range("C12").Value = 25.25677332204729
$ws.Range("D12").Value = 4.1009493093553
$ws.Range("E12").Value = 23.66167189386465
$ws.Range("F12").Value = 23.70072239388859
$ws.Range("G12").Value = 32.53780009353301
$ws.Range("H12").Value = 5.743300604138669
$ws.Range("J12").Value = 9.204087797426718
$ws.Range("K12").Value = 9.792174245106523
$ws.Range("B13").Value = 32.11880055635184
$ws.Range("C13").Value = 24.82079504961458
$ws.Range("D13").Value = 4.438850431706441
$ws.Range("E13").Value = 17.88552468053319
$ws.Range("F13").Value = 21.80955641832541
$ws.Range("G13").Value = 29.22140648844968
$ws.Range("H13").Value = 6.670476718896073
$ws.Range("J13").Value = 8.492959907930295
$ws.Range("K13").Value = 8.684601947464991
$ws.Range("B14").Value = 31.81564446389167
$ws.Range("C14").Value = 24.34939166932349
$ws.Range("D14").Value = 4.669126737831677
$ws.Range("E14").Value = 13.99216193476585
$ws.Range("F14").Value = 20.34472391562454
$ws.Range("G14").Value = 26.67600951845493
$ws.Range("H14").Value = 7.41299849926423
$ws.Range("J14").Value = 7.982689293514213
$ws.Range("K14").Value = 7.986287241427837
$ws.Range("B15").Value = 31.63771381202893
$ws.Range("C15").Value = 24.15526300413777
$ws.Range("D15").Value = 4.717446154945213
$ws.Range("E15").Value = 13.06492348169354
$ws.Range("F15").Value = 19.92953109308296
$ws.Range("G15").Value = 25.96274323596182
$ws.Range("H15").Value = 7.582678591472758
$ws.Range("J15").Value = 7.854617528057052
$ws.Range("K15").Value = 7.8461114016283
$ws.Range("B16").Value = 30.629395569307
$ws.Range("C16").Value = 23.40139339081595
$ws.Range("D16").Value = 4.618537742669315
$ws.Range("E16").Value = 12.83533761943203
$ws.Range("F16").Value = 19.46795634787076
$ws.Range("G16").Value = 25.24394199131528
$ws.Range("H16").Value = 7.283823204834468
$ws.Range("J16").Value = 7.841075559681389
$ws.Range("K16").Value = 8.093966636608263
$ws.Range("B17").Value = 30.00676511390768
$ws.Range("C17").Value = 23.07296721234827
$ws.Range("D17").Value = 4.415809250879967
$ws.Range("E17").Value = 14.74638761048396
$ws.Range("F17").Value = 19.93304125661172
$ws.Range("G17").Value = 26.12933050045496
$ws.Range("H17").Value = 6.606083638714226
$ws.Range("J17").Value = 8.106485770173501
$ws.Range("K17").Value = 8.634739994084343
$ws.Range("B18").Value = 29.66220706518938
$ws.Range("C18").Value = 23.08810240252694
$ws.Range("D18").Value = 4.111944657598617
$ws.Range("E18").Value = 19.11068570180034
$ws.Range("F18").Value = 21.24055357348675
$ws.Range("G18").Value = 28.49038982662054
$ws.Range("H18").Value = 5.623916112666137
$ws.Range("J18").Value = 8.638652163665419
$ws.Range("K18").Value = 9.513125510743548
$ws.Range("B19").Value = 29.57738784404856
$ws.Range("C19").Value = 23.38479406038057
$ws.Range("D19").Value = 3.769634517804396
$ws.Range("E19").Value = 25.27637941184491
$ws.Range("F19").Value = 23.10431173563382
$ws.Range("G19").Value = 31.81822060738086
$ws.Range("H19").Value = 4.742700942549884
$ws.Range("J19").Value = 9.346734362794038
$ws.Range("K19").Value = 10.6454611424501
$ws.Range("B20").Value = 30.22772543526278
$ws.Range("C20").Value = 24.44313044267467
$ws.Range("D20").Value = 3.341351980458814
$ws.Range("E20").Value = 35.4286402408185
$ws.Range("F20").Value = 26.37548755248248
$ws.Range("G20").Value = 37.55180922748339
$ws.Range("H20").Value = 4.587836909872584
$ws.Range("J20").Value = 10.51129414093993
$ws.Range("K20").Value = 12.42931672018553
$ws.Range("B21").Value = 32.0900919111083
$ws.Range("C21").Value = 25.98821067844236
$ws.Range("D21").Value = 3.331472572562832
$ws.Range("E21").Value = 38.38212853461187
$ws.Range("F21").Value = 28.05373080040202
$ws.Range("G21").Value = 40.25227777744853
$ws.Range("H21").Value = 5.01776559119018
$ws.Range("J21").Value = 10.93005843281387
$ws.Range("K21").Value = 12.73698605367121
$ws.Range("B22").Value = 33.25522987451089
$ws.Range("C22").Value = 26.93149019893336
$ws.Range("D22").Value = 3.343312342543743
$ws.Range("E22").Value = 39.8055828142779
$ws.Range("F22").Value = 29.01129639132077
$ws.Range("G22").Value = 41.78232383652254
$ws.Range("H22").Value = 5.273958164191376
$ws.Range("J22").Value = 11.16027336193374
$ws.Range("K22").Value = 12.87971353528034
$ws.Range("B23").Value = 32.63768297247393
$ws.Range("C23").Value = 26.43142936554
$ws.Range("D23").Value = 3.3370608997126
$ws.Range("E23").Value = 39.05048716797728
$ws.Range("F23").Value = 28.50059825009035
$ws.Range("G23").Value = 40.95891857236105
$ws.Range("H23").Value = 5.137462771915303
$ws.Range("J23").Value = 11.0366331472039
$ws.Range("K23").Value = 12.80082368064497
$ws.Range("B24").Value = 30.20334007044616
$ws.Range("C24").Value = 24.46248537645034
$ws.Range("D24").Value = 3.311907401058376
$ws.Range("E24").Value = 36.08779440998175
$ws.Range("F24").Value = 26.55819922399895
$ws.Range("G24").Value = 37.88048097787322
$ws.Range("H24").Value = 4.614808551451458
$ws.Range("J24").Value = 10.58504292760889
$ws.Range("K24").Value = 12.55975710032087
$ws.Range("B25").Value = 27.36864641100482
$ws.Range("C25").Value = 22.17437379896721
$ws.Range("D25").Value = 3.281698456767105
$ws.Range("E25").Value = 32.66514866645217
$ws.Range("F25").Value = 24.44688846950158
$ws.Range("G25").Value = 34.50881742531911
$ws.Range("H25").Value = 4.036249746741522
$ws.Range("J25").Value = 10.13213547738503
$ws.Range("K25").Value = 12.41051113838102

# Column I is no longer used in rows 2-25; clear its contents so the cells are removed
$ws.Range("I2:I25").ClearContents()

# New columns P and Q (rows 2-25) are all zero
$ws.Range("P2:P25").Value = 0
$ws.Range("Q2:Q25").Value = 0
